# Apply updated market-price derived values to the Leve profit sheets.
# Values correspond to a scheduled market-data refresh (no formulas are used;
# all cells in columns H:N hold static computed numbers).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2940.8948
$ws.Range("I106").Value = 2594.2222
$ws.Range("K106").Value = 2594.2222
$ws.Range("M106").Value = -1963.2222
$ws.Range("H114").Value = 40429
$ws.Range("J114").Value = 40429
$ws.Range("L114").Value = 40429
$ws.Range("N114").Value = -49107
$ws.Range("H132").Value = 4470181.5
$ws.Range("I132").Value = 5440528.5
$ws.Range("J132").Value = 6585.2
$ws.Range("K132").Value = 16321585.5
$ws.Range("L132").Value = 19755.6
$ws.Range("M132").Value = -16319055.5
$ws.Range("N132").Value = -24815.6
$ws.Range("H137").Value = 2557.4255
$ws.Range("I137").Value = 2584.5676
$ws.Range("K137").Value = 7753.702799999999
$ws.Range("M137").Value = -5203.702799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 918.14703
$ws.Range("I2").Value = 913.65216
$ws.Range("J2").Value = 927.5454999999999
$ws.Range("K2").Value = 913.65216
$ws.Range("L2").Value = 927.5454999999999
$ws.Range("M2").Value = -800.65216
$ws.Range("N2").Value = -1153.5455
$ws.Range("H53").Value = 26040.334
$ws.Range("I53").Value = 9039
$ws.Range("J53").Value = 60043
$ws.Range("K53").Value = 9039
$ws.Range("L53").Value = 60043
$ws.Range("M53").Value = -8357
$ws.Range("N53").Value = -61407
$ws.Range("H61").Value = 1466.5238
$ws.Range("I61").Value = 1415.9474
$ws.Range("J61").Value = 1947
$ws.Range("K61").Value = 1415.9474
$ws.Range("L61").Value = 1947
$ws.Range("M61").Value = -1203.9474
$ws.Range("N61").Value = -2371
$ws.Range("H116").Value = 918.14703
$ws.Range("I116").Value = 913.65216
$ws.Range("J116").Value = 927.5454999999999
$ws.Range("K116").Value = 913.65216
$ws.Range("L116").Value = 927.5454999999999
$ws.Range("M116").Value = 1380.34784
$ws.Range("N116").Value = -5515.5455
$ws.Range("H122").Value = 2719.2727
$ws.Range("I122").Value = 1801.5
$ws.Range("J122").Value = 5166.6665
$ws.Range("K122").Value = 5404.5
$ws.Range("L122").Value = 15499.9995
$ws.Range("M122").Value = -2954.5
$ws.Range("N122").Value = -20399.9995
$ws.Range("H136").Value = 1466.5238
$ws.Range("I136").Value = 1415.9474
$ws.Range("J136").Value = 1947
$ws.Range("K136").Value = 4247.8422
$ws.Range("L136").Value = 5841
$ws.Range("M136").Value = -1697.8422
$ws.Range("N136").Value = -10941

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 918.14703
$ws.Range("I3").Value = 913.65216
$ws.Range("J3").Value = 927.5454999999999
$ws.Range("K3").Value = 913.65216
$ws.Range("L3").Value = 927.5454999999999
$ws.Range("M3").Value = -799.65216
$ws.Range("N3").Value = -1155.5455
$ws.Range("H80").Value = 667.2
$ws.Range("I80").Value = 689.5
$ws.Range("J80").Value = 656.7059
$ws.Range("K80").Value = 689.5
$ws.Range("L80").Value = 656.7059
$ws.Range("M80").Value = 308.5
$ws.Range("N80").Value = -2652.7059
$ws.Range("H83").Value = 667.2
$ws.Range("I83").Value = 689.5
$ws.Range("J83").Value = 656.7059
$ws.Range("K83").Value = 3447.5
$ws.Range("L83").Value = 3283.5295
$ws.Range("M83").Value = 1544.5
$ws.Range("N83").Value = -13267.5295

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 36728.38
$ws.Range("I31").Value = 60316.06
$ws.Range("J31").Value = 3312.5
$ws.Range("K31").Value = 60316.06
$ws.Range("L31").Value = 3312.5
$ws.Range("M31").Value = -60021.06
$ws.Range("N31").Value = -3902.5
$ws.Range("H34").Value = 36728.38
$ws.Range("I34").Value = 60316.06
$ws.Range("J34").Value = 3312.5
$ws.Range("K34").Value = 60316.06
$ws.Range("L34").Value = 3312.5
$ws.Range("M34").Value = -60114.06
$ws.Range("N34").Value = -3716.5
$ws.Range("H58").Value = 1208.1818
$ws.Range("I58").Value = 1274.7941
$ws.Range("J58").Value = 981.7
$ws.Range("K58").Value = 1274.7941
$ws.Range("L58").Value = 981.7
$ws.Range("M58").Value = -1071.7941
$ws.Range("N58").Value = -1387.7
$ws.Range("H110").Value = 40702
$ws.Range("J110").Value = 40702
$ws.Range("L110").Value = 40702
$ws.Range("N110").Value = -48882
$ws.Range("H122").Value = 2737.8462
$ws.Range("I122").Value = 3011.2
$ws.Range("J122").Value = 1826.6666
$ws.Range("K122").Value = 9033.599999999999
$ws.Range("L122").Value = 5479.9998
$ws.Range("M122").Value = -6583.599999999999
$ws.Range("N122").Value = -10379.9998
$ws.Range("H134").Value = 6034.4287
$ws.Range("I134").Value = 1352.0667
$ws.Range("J134").Value = 17740.334
$ws.Range("K134").Value = 4056.2001
$ws.Range("L134").Value = 53221.00199999999
$ws.Range("M134").Value = -1521.2001
$ws.Range("N134").Value = -58291.00199999999
$ws.Range("H136").Value = 1208.1818
$ws.Range("I136").Value = 1274.7941
$ws.Range("J136").Value = 981.7
$ws.Range("K136").Value = 3824.3823
$ws.Range("L136").Value = 2945.1
$ws.Range("M136").Value = -1274.3823
$ws.Range("N136").Value = -8045.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1069207.4
$ws.Range("I131").Value = 480
$ws.Range("J131").Value = 1202798.4
$ws.Range("K131").Value = 1440
$ws.Range("L131").Value = 3608395.2
$ws.Range("M131").Value = 3600
$ws.Range("N131").Value = -3618475.2
$ws.Range("H132").Value = 1681.125
$ws.Range("I132").Value = 1149.6666
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 10346.9994
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -7816.999400000001
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5275.115
$ws.Range("I70").Value = 5152.8
$ws.Range("J70").Value = 5351.5625
$ws.Range("K70").Value = 5152.8
$ws.Range("L70").Value = 5351.5625
$ws.Range("M70").Value = -4882.8
$ws.Range("N70").Value = -5891.5625
$ws.Range("H73").Value = 5275.115
$ws.Range("I73").Value = 5152.8
$ws.Range("J73").Value = 5351.5625
$ws.Range("K73").Value = 5152.8
$ws.Range("L73").Value = 5351.5625
$ws.Range("M73").Value = -4216.8
$ws.Range("N73").Value = -7223.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 505.16666
$ws.Range("I46").Value = 550.5
$ws.Range("K46").Value = 550.5
$ws.Range("M46").Value = -362.5
$ws.Range("H60").Value = 100000
$ws.Range("J60").Value = 100000
$ws.Range("L60").Value = 100000
$ws.Range("N60").Value = -101018
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H132").Value = 4488.644
$ws.Range("I132").Value = 4919.439
$ws.Range("J132").Value = 3507.389
$ws.Range("K132").Value = 14758.317
$ws.Range("L132").Value = 10522.167
$ws.Range("M132").Value = -12228.317
$ws.Range("N132").Value = -15582.167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2437.972
$ws.Range("I132").Value = 2607.4167
$ws.Range("J132").Value = 1513.7273
$ws.Range("K132").Value = 7822.250100000001
$ws.Range("L132").Value = 4541.1819
$ws.Range("M132").Value = -5292.250100000001
$ws.Range("N132").Value = -9601.1819
$ws.Range("H140").Value = 29364.5
$ws.Range("J140").Value = 29364.5
$ws.Range("L140").Value = 29364.5
$ws.Range("N140").Value = -39724.5
$ws.Range("H141").Value = 73461.664
$ws.Range("J141").Value = 73461.664
$ws.Range("L141").Value = 73461.664
$ws.Range("N141").Value = -83821.664
